# Translate newdir/test.xlsx in es (add fi, se columns; fill in es/he_IL rows; add Greek/French helper columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before the current K column (es/he_IL), shifting
#    the existing K:L (es/he_IL) data over to M:N.
$ws.Range("K1:L1").EntireColumn.Insert()

# 2. Populate the two newly inserted header cells (fi / se) with the same
#    header style (s=1) as the other header cells, by copying the format
#    from an existing styled header cell.
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 11).Value2 = "fi"
$ws.Cells.Item(1, 12).Value2 = "se"

# 3. Remove the now-obsolete context columns B2 (test 1) and D2 (test 1 fr).
$ws.Cells.Item(2, 2).ClearFormats()
$ws.Cells.Item(2, 2).ClearContents()
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 4).ClearContents()

# 4. Row 2: update/clear the de_DE/el/jp columns and reset their styling to
#    the plain (unstyled) look used throughout the rest of the data rows.
#    H2 is already an empty inline string cell; just drop its formatting.
$ws.Cells.Item(2, 8).ClearFormats()

$ws.Cells.Item(2, 9).ClearFormats()
$ws.Cells.Item(2, 9).Value2 = "δοκιμασία"

# Leading apostrophe forces text (not a date) while typing, then clearing
# the formatting afterwards drops the quote-prefix styling it introduces.
$ws.Cells.Item(2, 10).Value2 = "'2023-11-09"
$ws.Cells.Item(2, 10).ClearFormats()

# es / he_IL values for row 2 (now columns M/N after the insert above)
$ws.Cells.Item(2, 13).ClearFormats()
$ws.Cells.Item(2, 13).Value2 = "examen"

$ws.Cells.Item(2, 14).ClearFormats()
$ws.Cells.Item(2, 14).Value2 = "בדיקה"

# 5. Row 3: fill in the de_DE/el/jp translations.
$ws.Cells.Item(3, 8).ClearFormats()
$ws.Cells.Item(3, 8).Value2 = "test 1 fr"

$ws.Cells.Item(3, 9).ClearFormats()
$ws.Cells.Item(3, 9).Value2 = "παράδειγμα"

$ws.Cells.Item(3, 10).Value2 = "'2023-11-09"
$ws.Cells.Item(3, 10).ClearFormats()

$ws.Cells.Item(3, 13).ClearFormats()
$ws.Cells.Item(3, 13).Value2 = "ejemplo"

$ws.Cells.Item(3, 14).ClearFormats()
$ws.Cells.Item(3, 14).Value2 = "דוגמא"

# 6. Row 4: fill in the de_DE/el/jp translations.
$ws.Cells.Item(4, 8).ClearFormats()
$ws.Cells.Item(4, 8).Value2 = "test 1 fr"

$ws.Cells.Item(4, 9).ClearFormats()
$ws.Cells.Item(4, 9).Value2 = "ψάρι"

$ws.Cells.Item(4, 10).Value2 = "'2023-11-09"
$ws.Cells.Item(4, 10).ClearFormats()

$ws.Cells.Item(4, 13).ClearFormats()
$ws.Cells.Item(4, 13).Value2 = "pez"

$ws.Cells.Item(4, 14).ClearFormats()
$ws.Cells.Item(4, 14).Value2 = "דג"
